# TOEFL Progress workbook update
# - Rename sheet "Progree" -> "Progress"
# - Switch active sheet/tab from "Progress" to "TPO Results"
# - Update several Progress-tracking cells on the "Progress" sheet
# - Update a couple of "Tasks they are related to" notes (shared strings)

$wb = $excel.ActiveWorkbook

$progress = $wb.Worksheets.Item(1)   # was named "Progree"
$tpo      = $wb.Worksheets.Item(2)   # "TPO Results"

# Fix the sheet name typo
$progress.Name = "Progress"

# --- Data updates on the Progress sheet ---

# Row 3: Oxford English Grammar Course CD
$progress.Range("C3").Value = 18
$progress.Range("D3").Value = 7.1

# Row 4: Oxford English Grammar Course Book
$progress.Range("D4").Value = 83

# Row 6: mark Writing as done
$progress.Range("I6").Value = "Y"

# Row 16: mark Reading as done
$progress.Range("F16").Value = "Y"

# Row 17: mark Listening as done
$progress.Range("G17").Value = "Y"

# Row 18: mark Speaking as done
$progress.Range("H18").Value = "Y"

# Row 19: Reading/Listening no longer done, Writing now done
$progress.Range("F19").ClearContents()
$progress.Range("G19").ClearContents()
$progress.Range("I19").Value = "Y"

# Row 23: mark Writing as done
$progress.Range("I23").Value = "Y"

# Row 24: progress update
$progress.Range("D24").Value = 200

# Row 37: mark Speaking and Writing as done
# (set H37 before I37 so new shared strings keep the same relative order
# as the updated "tasks" notes below)
$progress.Range("H37").Value = "Y"
$progress.Range("I37").Value = "Y"

# Update the "Tasks they are related to" notes (J4 before J3 keeps the
# shared-string table ordering consistent with the source workbook)
$progress.Range("J4").Value = "17, 28, 29, 34, 36, 60, 74"
$progress.Range("J3").Value = "1.5, 2.7, 3.7, 6.1, 6.7"

# --- Window / selection state ---
# Leave a selection behind on the Progress sheet, then switch to and
# select a cell on the TPO Results sheet, which becomes the active tab.
$progress.Range("H11").Select()
$tpo.Activate()
$tpo.Range("C11").Select()
